$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2021.1072
$ws.Range("I15").Value = 2021.1072
$ws.Range("K15").Value = 6063.321599999999
$ws.Range("M15").Value = -5894.321599999999
# Row 86
$ws.Range("H86").Value = 2337.625
$ws.Range("I86").Value = 2243
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2243
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1120
$ws.Range("N86").Value = -5246
# Row 87
$ws.Range("H87").Value = 22910.842
$ws.Range("J87").Value = 22910.842
$ws.Range("L87").Value = 22910.842
$ws.Range("N87").Value = -25406.842
# Row 89
$ws.Range("H89").Value = 2337.625
$ws.Range("I89").Value = 2243
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 11215
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -5599
$ws.Range("N89").Value = -26232
# Row 90
$ws.Range("H90").Value = 22910.842
$ws.Range("J90").Value = 22910.842
$ws.Range("L90").Value = 68732.526
$ws.Range("N90").Value = -81212.526
# Row 103
$ws.Range("H103").Value = 419.7
$ws.Range("I103").Value = 264
$ws.Range("J103").Value = 471.6
$ws.Range("K103").Value = 792
$ws.Range("L103").Value = 1414.8
$ws.Range("M103").Value = -206
$ws.Range("N103").Value = -2586.8
# Row 106
$ws.Range("H106").Value = 2900
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 2900
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 2900
$ws.Range("N106").Value = -4162
$ws.Range("M106").ClearContents()
# Row 111
$ws.Range("H111").Value = 627.55
$ws.Range("J111").Value = 787.4286
$ws.Range("L111").Value = 2362.2858
$ws.Range("N111").Value = -8496.2858
# Row 137
$ws.Range("H137").Value = 4421.242
$ws.Range("I137").Value = 4496.1924
$ws.Range("J137").Value = 4142.857
$ws.Range("K137").Value = 13488.5772
$ws.Range("L137").Value = 12428.571
$ws.Range("M137").Value = -10938.5772
$ws.Range("N137").Value = -17528.571

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 10870
$ws.Range("I61").Value = 19435.334
$ws.Range("J61").Value = 2304.6667
$ws.Range("K61").Value = 19435.334
$ws.Range("L61").Value = 2304.6667
$ws.Range("M61").Value = -19223.334
$ws.Range("N61").Value = -2728.6667
# Row 97
$ws.Range("H97").Value = 1269.7142
$ws.Range("I97").Value = 1269.7142
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1269.7142
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -773.7141999999999
$ws.Range("N97").ClearContents()
# Row 112
$ws.Range("H112").Value = 28387
$ws.Range("J112").Value = 28387
$ws.Range("L112").Value = 28387
$ws.Range("N112").Value = -31341
# Row 136
$ws.Range("H136").Value = 10870
$ws.Range("I136").Value = 19435.334
$ws.Range("J136").Value = 2304.6667
$ws.Range("K136").Value = 58306.00199999999
$ws.Range("L136").Value = 6914.000100000001
$ws.Range("M136").Value = -55756.00199999999
$ws.Range("N136").Value = -12014.0001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 686.2105
$ws.Range("J80").Value = 594.0909
$ws.Range("L80").Value = 594.0909
$ws.Range("N80").Value = -2590.0909
# Row 82
$ws.Range("H82").Value = 25300.857
$ws.Range("J82").Value = 28641.5
$ws.Range("L82").Value = 28641.5
$ws.Range("N82").Value = -29407.5
# Row 83
$ws.Range("H83").Value = 686.2105
$ws.Range("J83").Value = 594.0909
$ws.Range("L83").Value = 2970.4545
$ws.Range("N83").Value = -12954.4545
# Row 85
$ws.Range("H85").Value = 25300.857
$ws.Range("J85").Value = 28641.5
$ws.Range("L85").Value = 28641.5
$ws.Range("N85").Value = -31293.5
# Row 86
$ws.Range("H86").Value = 28963.525
$ws.Range("I86").Value = 2806.9285
$ws.Range("J86").Value = 102202
$ws.Range("K86").Value = 2806.9285
$ws.Range("L86").Value = 102202
$ws.Range("M86").Value = -1683.9285
$ws.Range("N86").Value = -104448
# Row 89
$ws.Range("H89").Value = 28963.525
$ws.Range("I89").Value = 2806.9285
$ws.Range("J89").Value = 102202
$ws.Range("K89").Value = 14034.6425
$ws.Range("L89").Value = 511010
$ws.Range("M89").Value = -8418.6425
$ws.Range("N89").Value = -522242
# Row 94
$ws.Range("H94").Value = 917.5
$ws.Range("I94").Value = 590.5
$ws.Range("J94").Value = 1244.5
$ws.Range("K94").Value = 590.5
$ws.Range("L94").Value = 1244.5
$ws.Range("M94").Value = -139.5
$ws.Range("N94").Value = -2146.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 248.85715
$ws.Range("J98").Value = 348.57144
$ws.Range("L98").Value = 1045.71432
$ws.Range("N98").Value = -4041.71432
# Row 111
$ws.Range("H111").Value = 1685
$ws.Range("I111").Value = 296
$ws.Range("K111").Value = 888
$ws.Range("M111").Value = 2179
# Row 114
$ws.Range("H114").Value = 765.35297
$ws.Range("J114").Value = 1805.1666
$ws.Range("L114").Value = 5415.4998
$ws.Range("N114").Value = -11923.4998

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 28251.25
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2712
# Row 58
$ws.Range("H58").Value = 80046
$ws.Range("J58").Value = 80046
$ws.Range("L58").Value = 80046
$ws.Range("N58").Value = -80600
# Row 80
$ws.Range("H80").Value = 3995
$ws.Range("J80").Value = 4990
$ws.Range("L80").Value = 4990
$ws.Range("N80").Value = -6986
# Row 83
$ws.Range("H83").Value = 3995
$ws.Range("J83").Value = 4990
$ws.Range("L83").Value = 24950
$ws.Range("N83").Value = -34934
# Row 113
$ws.Range("H113").Value = 1945.8948
$ws.Range("I113").Value = 2029.5
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2029.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 140.5
$ws.Range("N113").Value = -5840

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 324.5
$ws.Range("I22").Value = 324.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 324.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -29.5
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 324.5
$ws.Range("I27").Value = 324.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 324.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -217.5
$ws.Range("N27").ClearContents()
# Row 68
$ws.Range("H68").Value = 1581.5834
$ws.Range("I68").Value = 1075.3636
$ws.Range("K68").Value = 1075.3636
$ws.Range("M68").Value = -326.3635999999999
# Row 71
$ws.Range("H71").Value = 1581.5834
$ws.Range("I71").Value = 1075.3636
$ws.Range("K71").Value = 5376.817999999999
$ws.Range("M71").Value = -1632.817999999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 1680502.6
$ws.Range("I3").Value = 5000050
$ws.Range("J3").Value = 20729
$ws.Range("K3").Value = 5000050
$ws.Range("L3").Value = 20729
$ws.Range("M3").Value = -4999936
$ws.Range("N3").Value = -20957
# Row 62
$ws.Range("H62").Value = 5068.1665
$ws.Range("I62").Value = 5500
$ws.Range("K62").Value = 5500
$ws.Range("M62").Value = -4876
# Row 65
$ws.Range("H65").Value = 5068.1665
$ws.Range("I65").Value = 5500
$ws.Range("K65").Value = 27500
$ws.Range("M65").Value = -24380
# Row 136
$ws.Range("H136").Value = 3703.0667
$ws.Range("I136").Value = 3413.2727
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 10239.8181
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -7689.8181
$ws.Range("N136").Value = -18600
